# Nueva seccion multiproductos y cambios en el escaner de codigos
# Adds 19 new "DOBLE A" aerosol products as rows 7-25 on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Codigo, Proveedor, Nombre, Precio, Observaciones, Dueno, Cantidad
$products = @(
    @("5436", "DOBLE A", "Aerosol DOBLEA COM - NEGRO BRILL",   "5,500", "", "ferreteria_general", 1),
    @("5435", "DOBLE A", "Aerosol DOBLEA COM - NEGRO MATE",    "5,500", "", "ferreteria_general", 1),
    @("0026", "DOBLE A", "Aerosol DOBLEA COM - NEGRO SATIN",   "5,500", "", "ferreteria_general", 1),
    @("5437", "DOBLE A", "Aerosol DOBLEA COM - BLANCO BRILL",  "5,500", "", "ferreteria_general", 1),
    @("5445", "DOBLE A", "Aerosol DOBLEA COM - BLANCO MATE",   "5,500", "", "ferreteria_general", 1),
    @("0028", "DOBLE A", "Aerosol DOBLEA COM - BLANCO SATIN",  "5,500", "", "ferreteria_general", 1),
    @("5441", "DOBLE A", "Aerosol DOBLEA COM - AMARILLO",      "5,500", "", "ferreteria_general", 1),
    @("5451", "DOBLE A", "Aerosol DOBLEA COM - VERDE",         "5,500", "", "ferreteria_general", 1),
    @("5443", "DOBLE A", "Aerosol DOBLEA COM - AZUL",          "5,500", "", "ferreteria_general", 1),
    @("5439", "DOBLE A", "Aerosol DOBLEA COM - ROJO",          "5,500", "", "ferreteria_general", 1),
    @("5452", "DOBLE A", "Aerosol DOBLEA COM - BEIGE",         "5,500", "", "ferreteria_general", 1),
    @("5444", "DOBLE A", "Aerosol DOBLEA COM - AZUL MARINO",   "5,500", "", "ferreteria_general", 1),
    @("5450", "DOBLE A", "Aerosol DOBLEA COM - VERDE OSCURO",  "5,500", "", "ferreteria_general", 1),
    @("5453", "DOBLE A", "Aerosol DOBLEA COM - MARRON",        "5,500", "", "ferreteria_general", 1),
    @("5447", "DOBLE A", "Aerosol DOBLEA COM - GRAFITO",       "5,500", "", "ferreteria_general", 1),
    @("5448", "DOBLE A", "Aerosol DOBLEA COM - NARANJA",       "5,500", "", "ferreteria_general", 1),
    @("5449", "DOBLE A", "Aerosol DOBLEA COM - ROSA",          "5,500", "", "ferreteria_general", 1),
    @("5438", "DOBLE A", "Aerosol DOBLEA COM - ALUM METAL",    "5,500", "", "ferreteria_general", 1),
    @("5440", "DOBLE A", "Aerosol DOBLEA COM - DORADO",        "5,500", "", "ferreteria_general", 1)
)

$startRow = 7
$endRow = $startRow + $products.Length - 1

# Force columns A, D and E to be treated as plain text so values such as
# "0026" keep their leading zeros and "5,500" is not parsed as a number.
$ws.Range("A$startRow" + ":A$endRow").NumberFormat = "@"
$ws.Range("D$startRow" + ":D$endRow").NumberFormat = "@"
$ws.Range("E$startRow" + ":E$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $products.Length; $i++) {
    $row = $startRow + $i
    $item = $products[$i]

    $ws.Cells.Item($row, 1).Value = $item[0]   # Codigo
    $ws.Cells.Item($row, 2).Value = $item[1]   # Proveedor
    $ws.Cells.Item($row, 3).Value = $item[2]   # Nombre
    $ws.Cells.Item($row, 4).Value = $item[3]   # Precio
    $ws.Cells.Item($row, 5).Value = $item[4]   # Observaciones
    $ws.Cells.Item($row, 6).Value = $item[5]   # Dueno
    $ws.Cells.Item($row, 7).Value = $item[6]   # Cantidad (numero)
}

# Restore default styling on the text-forced columns so the new cells don't
# carry a stray custom number format, matching the original sheet's look.
$ws.Range("A$startRow" + ":A$endRow").Style = "Normal"
$ws.Range("D$startRow" + ":D$endRow").Style = "Normal"
$ws.Range("E$startRow" + ":E$endRow").Style = "Normal"
